# "Update countries & provincias Spain" -- refresh the COVID country
# dashboard on sheet "Pais": a handful of countries swap rank (so their
# name cells trade places) and the daily case/recovered/death counters
# move forward; the "last updated" timestamp also advances.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 30 de Septiembre de 2020 a las 18:40"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7417308
$ws.Range("C4").Value = 11162
$ws.Range("D4").Value = 4653343
$ws.Range("E4").Value = 2552783
$ws.Range("G4").Value = 397
$ws.Range("H4").Value = 211182

# Row 5 - India
$ws.Range("B5").Value = 6269042
$ws.Range("C5").Value = 45523
$ws.Range("D5").Value = 5225766
$ws.Range("E5").Value = 945333
$ws.Range("G5").Value = 414
$ws.Range("H5").Value = 97943

# Row 6 - Brasil
$ws.Range("B6").Value = 4787637
$ws.Range("C6").Value = 7320
$ws.Range("E6").Value = 509306
$ws.Range("G6").Value = 233
$ws.Range("H6").Value = 143243

# Row 17 - Reino Unido
$ws.Range("B17").Value = 453264
$ws.Range("C17").Value = 7108
$ws.Range("G17").Value = 71
$ws.Range("H17").Value = 42143

# Row 25 - Alemania
$ws.Range("B25").Value = 291578
$ws.Range("C25").Value = 1112
$ws.Range("E25").Value = 26019

# Row 29 - Canada
$ws.Range("B29").Value = 158425
$ws.Range("C29").Value = 1464
$ws.Range("E29").Value = 14936

# Row 30 - Ecuador
$ws.Range("B30").Value = 137047
$ws.Range("C30").Value = 1298
$ws.Range("E30").Value = 13396
$ws.Range("G30").Value = 43
$ws.Range("H30").Value = 11355

# Row 57 - Chequia
$ws.Range("B57").Value = 68919
$ws.Range("C57").Value = 1076
$ws.Range("D57").Value = 33443
$ws.Range("E57").Value = 34821
$ws.Range("G57").Value = 19
$ws.Range("H57").Value = 655

# Row 60 - Uzbekistan
$ws.Range("B60").Value = 56717
$ws.Range("C60").Value = 363
$ws.Range("D60").Value = 53366
$ws.Range("E60").Value = 2881
$ws.Range("G60").Value = 4
$ws.Range("H60").Value = 470

# Row 86 - Grecia
$ws.Range("B86").Value = 18475
$ws.Range("C86").Value = 352
$ws.Range("E86").Value = 8095
$ws.Range("G86").Value = 3
$ws.Range("H86").Value = 391

# Rows 100-102 swap rank: Montenegro moves above Consejo Danes / Guinea
$ws.Range("A100").Value = "Montenegro"
$ws.Range("B100").Value = 10772
$ws.Range("C100").Value = 197
$ws.Range("D100").Value = 7192
$ws.Range("E100").Value = 3411
$ws.Range("G100").Value = 5
$ws.Range("H100").Value = 169

$ws.Range("A101").Value = "Consejo Danes para los Refugiados"
$ws.Range("B101").Value = 10659
$ws.Range("C101").Value = 28
$ws.Range("D101").Value = 10139
$ws.Range("E101").Value = 248
$ws.Range("H101").Value = 272

$ws.Range("A102").Value = "Guinea"
$ws.Range("B102").Value = 10634
$ws.Range("D102").Value = 9960
$ws.Range("E102").Value = 608
$ws.Range("H102").Value = 66

# Row 106 - Guayana Francesa
$ws.Range("B106").Value = 9955
$ws.Range("C106").Value = 26
$ws.Range("D106").Value = 9589
$ws.Range("E106").Value = 300

# Row 111 - Luxemburgo
$ws.Range("B111").Value = 8509
$ws.Range("C111").Value = 78
$ws.Range("D111").Value = 7174
$ws.Range("E111").Value = 1211

# Row 118 - Malaui
$ws.Range("B118").Value = 5773
$ws.Range("C118").Value = 1
$ws.Range("D118").Value = 4263
$ws.Range("E118").Value = 1331

# Rows 138-141 swap rank: Mayotte moves above Somalia / Gambia / Tailandia
$ws.Range("A138").Value = "Mayotte"
$ws.Range("B138").Value = 3779
$ws.Range("C138").Value = 238
$ws.Range("D138").Value = 2964
$ws.Range("E138").Value = 773
$ws.Range("G138").Value = 2
$ws.Range("H138").Value = 42

$ws.Range("A139").Value = "Somalia"
$ws.Range("B139").Value = 3588
$ws.Range("D139").Value = 2946
$ws.Range("E139").Value = 543
$ws.Range("H139").Value = 99

$ws.Range("A140").Value = "Gambia"
$ws.Range("B140").Value = 3579
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 2161
$ws.Range("E140").Value = 1306
$ws.Range("H140").Value = 112

$ws.Range("A141").Value = "Tailandia"
$ws.Range("B141").Value = 3564
$ws.Range("C141").Value = 5
$ws.Range("D141").Value = 3374
$ws.Range("E141").Value = 131
$ws.Range("H141").Value = 59

# Row 142 - Sri Lanka
$ws.Range("B142").Value = 3380
$ws.Range("C142").Value = 6
$ws.Range("E142").Value = 137

# Rows 153-155 swap rank: Yemen moves above Uruguay / Burkina Faso
$ws.Range("A153").Value = "Yemen"
$ws.Range("B153").Value = 2034
$ws.Range("C153").Value = 3
$ws.Range("D153").Value = 1286
$ws.Range("E153").Value = 161
$ws.Range("H153").Value = 587

$ws.Range("A154").Value = "Uruguay"
$ws.Range("B154").Value = 2033
$ws.Range("D154").Value = 1771
$ws.Range("E154").Value = 214
$ws.Range("H154").Value = 48

$ws.Range("A155").Value = "Burkina Faso"
$ws.Range("B155").Value = 2032
$ws.Range("D155").Value = 1309
$ws.Range("E155").Value = 665
$ws.Range("H155").Value = 58

# Row 161 - Belice
$ws.Range("B161").Value = 1755
$ws.Range("C161").Value = 12
$ws.Range("E161").Value = 364

# Row 163 - Republica de Chipre
$ws.Range("B163").Value = 1595
$ws.Range("C163").Value = 19
$ws.Range("E163").Value = 686
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 36

# Rows 164-165 swap rank: Martinica moves above Liberia
$ws.Range("A164").Value = "Martinica"
$ws.Range("B164").Value = 1543
$ws.Range("C164").Value = 253
$ws.Range("D164").Value = 98
$ws.Range("E164").Value = 1424
$ws.Range("G164").Value = 1
$ws.Range("H164").Value = 21

$ws.Range("A165").Value = "Liberia"
$ws.Range("B165").Value = 1343
$ws.Range("D165").Value = 1221
$ws.Range("E165").Value = 40
$ws.Range("H165").Value = 82
